# RediffTestData.xlsx - Day1 test-data refresh:
#  - swap out the "seleniumpanda2@rediffmail.com" row for a new account
#    (jsmith2024@rediffmail.com / redcow@1999)
#  - append a brand-new credential row (homa_rahimi@rediffmail.com / Winter2021?)
#    with a mailto hyperlink on the e-mail cell, mirroring the existing rows
#  - leave the recorded UI state (selection, page orientation) as captured

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: replace the retired seleniumpanda2 account ---------------------
$ws.Range("A4").Value = "jsmith2024@rediffmail.com"
$ws.Range("B4").Value = "redcow@1999"

# --- Row 6 (new): append another login pair --------------------------------
$ws.Range("A6").Value = "homa_rahimi@rediffmail.com"
$ws.Range("B6").Value = "Winter2021?"

# Hyperlink the new e-mail cell just like the other username/password cells
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:homa_rahimi@rediffmail.com")

# --- Misc UI/page state recorded by the workbook author ---------------------
$ws.Range("B18").Select()
$ws.PageSetup.Orientation = 1
